$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

# Columns A-D hold text that looks like dates/numbers ("2023-06-30", "26", ...).
# Force them to be written as literal text (not auto-converted to a date serial
# or a number) by temporarily switching the cell to a text number format, then
# clearing the format again afterwards so the cell keeps the default style.
$textCols = 1,2,3,4
$textVals = @("2023-06-30","08:32:05","Friday","26")
for ($i = 0; $i -lt $textCols.Length; $i++) {
    $cell = $ws.Cells.Item($row, $textCols[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $textVals[$i]
    $cell.ClearFormats()
}

# Columns E-T are plain numbers.
$ws.Cells.Item($row, 5).Value = 123313
$ws.Cells.Item($row, 6).Value = 134468
$ws.Cells.Item($row, 7).Value = 163616
$ws.Cells.Item($row, 8).Value = 133976
$ws.Cells.Item($row, 9).Value = 177029
$ws.Cells.Item($row, 10).Value = 115417
$ws.Cells.Item($row, 11).Value = 204560
$ws.Cells.Item($row, 12).Value = 226104
$ws.Cells.Item($row, 13).Value = 176701
$ws.Cells.Item($row, 14).Value = 104611
$ws.Cells.Item($row, 15).Value = 39806
$ws.Cells.Item($row, 16).Value = 33703
$ws.Cells.Item($row, 17).Value = 52580
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35703
$ws.Cells.Item($row, 20).Value = -1
